{"js": "// Old/new text for each \"AA\u00d7BB=\" multiplication cell, taken from the\n// authoritative commit diff. Every \"old\" string is unique in the document\n// and none of the \"new\" strings collide with any other \"old\" string, so a\n// simple sequential whole-text search/replace per pair is safe regardless\n// of processing order.\nconst pairs = [\n  [\"82\u00d793=\", \"31\u00d789=\"],\n  [\"96\u00d717=\", \"35\u00d794=\"],\n  [\"65\u00d732=\", \"81\u00d717=\"],\n  [\"82\u00d763=\", \"55\u00d716=\"],\n  [\"12\u00d760=\", \"84\u00d777=\"],\n  [\"11\u00d718=\", \"23\u00d779=\"],\n  [\"39\u00d786=\", \"48\u00d726=\"],\n  [\"57\u00d739=\", \"15\u00d796=\"],\n  [\"55\u00d713=\", \"66\u00d759=\"],\n  [\"23\u00d752=\", \"17\u00d739=\"],\n  [\"98\u00d713=\", \"97\u00d737=\"],\n  [\"56\u00d736=\", \"99\u00d729=\"],\n  [\"49\u00d764=\", \"84\u00d716=\"],\n  [\"38\u00d770=\", \"85\u00d779=\"],\n  [\"34\u00d762=\", \"59\u00d781=\"],\n  [\"82\u00d777=\", \"36\u00d774=\"],\n  [\"62\u00d781=\", \"40\u00d733=\"],\n  [\"42\u00d799=\", \"67\u00d719=\"],\n  [\"81\u00d740=\", \"72\u00d736=\"],\n  [\"25\u00d751=\", \"90\u00d782=\"],\n  [\"59\u00d751=\", \"86\u00d731=\"],\n  [\"91\u00d721=\", \"19\u00d715=\"],\n  [\"17\u00d766=\", \"89\u00d738=\"],\n  [\"93\u00d788=\", \"96\u00d795=\"],\n  [\"98\u00d769=\", \"13\u00d723=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Old/new pairs for each \"AA\u00d7BB=\" cell in the practice-sheet table, taken\n# from the authoritative commit diff. Every \"old\" string is unique in the\n# document, so a plain whole-content Find/Replace per pair is safe and\n# order-independent (no new value collides with any other old value).\n$replacements = @(\n    @(\"82\u00d793=\", \"31\u00d789=\"),\n    @(\"96\u00d717=\", \"35\u00d794=\"),\n    @(\"65\u00d732=\", \"81\u00d717=\"),\n    @(\"82\u00d763=\", \"55\u00d716=\"),\n    @(\"12\u00d760=\", \"84\u00d777=\"),\n    @(\"11\u00d718=\", \"23\u00d779=\"),\n    @(\"39\u00d786=\", \"48\u00d726=\"),\n    @(\"57\u00d739=\", \"15\u00d796=\"),\n    @(\"55\u00d713=\", \"66\u00d759=\"),\n    @(\"23\u00d752=\", \"17\u00d739=\"),\n    @(\"98\u00d713=\", \"97\u00d737=\"),\n    @(\"56\u00d736=\", \"99\u00d729=\"),\n    @(\"49\u00d764=\", \"84\u00d716=\"),\n    @(\"38\u00d770=\", \"85\u00d779=\"),\n    @(\"34\u00d762=\", \"59\u00d781=\"),\n    @(\"82\u00d777=\", \"36\u00d774=\"),\n    @(\"62\u00d781=\", \"40\u00d733=\"),\n    @(\"42\u00d799=\", \"67\u00d719=\"),\n    @(\"81\u00d740=\", \"72\u00d736=\"),\n    @(\"25\u00d751=\", \"90\u00d782=\"),\n    @(\"59\u00d751=\", \"86\u00d731=\"),\n    @(\"91\u00d721=\", \"19\u00d715=\"),\n    @(\"17\u00d766=\", \"89\u00d738=\"),\n    @(\"93\u00d788=\", \"96\u00d795=\"),\n    @(\"98\u00d769=\", \"13\u00d723=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
